$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.645.32'
$ws.Range("E2").Value = '  -0.07%  '

$ws.Range("D3").Value = '1.642.31'
$ws.Range("E3").Value = '  +0.68%  '

$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.71'
$ws.Range("E5").Value = '  +0.63%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.504'
$ws.Range("E6").Value = '  +0.82%  '

$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("E9").Value = '  +0.76%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.23'
$ws.Range("E10").Value = '  +0.10%  '

$ws.Range("E11").Value = '  +0.19%  '

$ws.Range("D12").Value = '1.871.50'

$ws.Range("E13").Value = '  +2.71%  '

$ws.Range("D14").Value = '1.636.42'
$ws.Range("E14").Value = '  +0.83%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.529'
$ws.Range("E15").Value = '  +1.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.48'
$ws.Range("E16").Value = '  +3.17%  '

$ws.Range("D17").Value = '26.693.65'
$ws.Range("E17").Value = '  +0.15%  '

$ws.Range("E18").Value = '  +0.98%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.80'
$ws.Range("E19").Value = '  -1.30%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.00'
$ws.Range("E20").Value = '  +0.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.36'
$ws.Range("E21").Value = '  +1.59%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.29'
$ws.Range("E22").Value = '  +2.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.54'
$ws.Range("E23").Value = '  +1.87%  '

$ws.Range("E24").Value = '  +11.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.91'
$ws.Range("E25").Value = '  -1.29%  '

$ws.Range("E26").Value = '  +0.27%  '

$ws.Range("E27").Value = '  -0.50%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.76'
$ws.Range("E29").Value = '  +1.46%  '

$ws.Range("E30").Value = '  +2.42%  '

$ws.Range("E31").Value = '  +0.22%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.38'
$ws.Range("E32").Value = '  +2.56%  '

$ws.Range("E33").Value = '  +2.17%  '

$ws.Range("D34").Value = '1.274.91'
$ws.Range("E34").Value = '  +4.97%  '

$ws.Range("E35").Value = '  +2.08%  '

$ws.Range("E36").Value = '  +5.79%  '

$ws.Range("E37").Value = '  +0.16%  '

$ws.Range("E38").Value = '  +6.14%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.828'
$ws.Range("E39").Value = '  +2.93%  '

$ws.Range("E40").Value = '  +0.21%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.814'
$ws.Range("E41").Value = '  +2.49%  '

$ws.Range("E42").Value = '  -1.70%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.45'
$ws.Range("E43").Value = '  +2.02%  '

$ws.Range("D44").Value = '1.781.98'
$ws.Range("E44").Value = '  +0.78%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.71'
$ws.Range("E45").Value = '  -0.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.47'
$ws.Range("E46").Value = '  +7.98%  '

$ws.Range("E47").Value = '  +2.51%  '

$ws.Range("E48").Value = '  +0.81%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.84'
$ws.Range("E49").Value = '  +2.96%  '

$ws.Range("E50").Value = '  +3.04%  '

$ws.Range("E51").Value = '  -0.65%  '
